$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nid2"
$ws.Range("C2").Value = "Col13a1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.88772733333334
$ws.Range("H2").Value = 140.663182
$ws.Range("I2").Value = 0.4281561666633809
$ws.Range("J2").Value = 0.4281561666633809
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.335983
$ws.Range("N2").Value = 1.007949
$ws.Range("O2").Value = 0.4945457382278176
$ws.Range("P2").Value = 0.4945457382278176
$ws.Range("Q2").Value = 15.75347929263533
$ws.Range("R2").Value = 141.781313633718
$ws.Range("S2").Value = 0.2117428075193342
$ws.Range("T2").Value = 0.2117428075193342

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nid2"
$ws.Range("C3").Value = "Col13a1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.88772733333334
$ws.Range("H3").Value = 140.663182
$ws.Range("I3").Value = 0.4281561666633809
$ws.Range("J3").Value = 0.4281561666633809
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.343394
$ws.Range("N3").Value = 1.030182
$ws.Range("O3").Value = 0.5054542617721824
$ws.Range("P3").Value = 0.5054542617721824
$ws.Range("Q3").Value = 16.10096423990267
$ws.Range("R3").Value = 144.908678159124
$ws.Range("S3").Value = 0.2164133591440467
$ws.Range("T3").Value = 0.2164133591440467

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nid2"
$ws.Range("C4").Value = "Col13a1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 54.562069
$ws.Range("H4").Value = 163.686207
$ws.Range("I4").Value = 0.4982345623660686
$ws.Range("J4").Value = 0.4982345623660687
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.335983
$ws.Range("N4").Value = 1.007949
$ws.Range("O4").Value = 0.4945457382278176
$ws.Range("P4").Value = 0.4945457382278176
$ws.Range("Q4").Value = 18.331927628827
$ws.Range("R4").Value = 164.987348659443
$ws.Range("S4").Value = 0.246399779455941
$ws.Range("T4").Value = 0.246399779455941

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nid2"
$ws.Range("C5").Value = "Col13a1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 54.562069
$ws.Range("H5").Value = 163.686207
$ws.Range("I5").Value = 0.4982345623660686
$ws.Range("J5").Value = 0.4982345623660687
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.343394
$ws.Range("N5").Value = 1.030182
$ws.Range("O5").Value = 0.5054542617721824
$ws.Range("P5").Value = 0.5054542617721824
$ws.Range("Q5").Value = 18.736287122186
$ws.Range("R5").Value = 168.626584099674
$ws.Range("S5").Value = 0.2518347829101276
$ws.Range("T5").Value = 0.2518347829101276

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Nid2"
$ws.Range("C6").Value = "Col13a1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7240593333333334
$ws.Range("H6").Value = 2.172178
$ws.Range("I6").Value = 0.006611761461435795
$ws.Range("J6").Value = 0.006611761461435797
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.335983
$ws.Range("N6").Value = 1.007949
$ws.Range("O6").Value = 0.4945457382278176
$ws.Range("P6").Value = 0.4945457382278176
$ws.Range("Q6").Value = 0.2432716269913333
$ws.Range("R6").Value = 2.189444642922
$ws.Range("S6").Value = 0.003269818452931999
$ws.Range("T6").Value = 0.003269818452932

$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Nid2"
$ws.Range("C7").Value = "Col13a1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7240593333333334
$ws.Range("H7").Value = 2.172178
$ws.Range("I7").Value = 0.006611761461435795
$ws.Range("J7").Value = 0.006611761461435797
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.343394
$ws.Range("N7").Value = 1.030182
$ws.Range("O7").Value = 0.5054542617721824
$ws.Range("P7").Value = 0.5054542617721824
$ws.Range("Q7").Value = 0.2486376307106667
$ws.Range("R7").Value = 2.237738676396
$ws.Range("S7").Value = 0.003341943008503796
$ws.Range("T7").Value = 0.003341943008503797

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Nid2"
$ws.Range("C8").Value = "Col13a1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.037292333333333
$ws.Range("H8").Value = 3.111877
$ws.Range("I8").Value = 0.009472054509956569
$ws.Range("J8").Value = 0.00947205450995657
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.335983
$ws.Range("N8").Value = 1.007949
$ws.Range("O8").Value = 0.4945457382278176
$ws.Range("P8").Value = 0.4945457382278176
$ws.Range("Q8").Value = 0.3485125900303334
$ws.Range("R8").Value = 3.136613310273
$ws.Range("S8").Value = 0.004684364190160601
$ws.Range("T8").Value = 0.004684364190160601

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Nid2"
$ws.Range("C9").Value = "Col13a1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.037292333333333
$ws.Range("H9").Value = 3.111877
$ws.Range("I9").Value = 0.009472054509956569
$ws.Range("J9").Value = 0.00947205450995657
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.343394
$ws.Range("N9").Value = 1.030182
$ws.Range("O9").Value = 0.5054542617721824
$ws.Range("P9").Value = 0.5054542617721824
$ws.Range("Q9").Value = 0.3561999635126667
$ws.Range("R9").Value = 3.205799671614
$ws.Range("S9").Value = 0.004787690319795969
$ws.Range("T9").Value = 0.00478769031979597

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Nid2"
$ws.Range("C10").Value = "Col13a1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.299658999999999
$ws.Range("H10").Value = 18.898977
$ws.Range("I10").Value = 0.0575254549991582
$ws.Range("J10").Value = 0.05752545499915822
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.335983
$ws.Range("N10").Value = 1.007949
$ws.Range("O10").Value = 0.4945457382278176
$ws.Range("P10").Value = 0.4945457382278176
$ws.Range("Q10").Value = 2.116578329797
$ws.Range("R10").Value = 19.049204968173
$ws.Range("S10").Value = 0.02844896860944979
$ws.Range("T10").Value = 0.0284489686094498

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Nid2"
$ws.Range("C11").Value = "Col13a1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.299658999999999
$ws.Range("H11").Value = 18.898977
$ws.Range("I11").Value = 0.0575254549991582
$ws.Range("J11").Value = 0.05752545499915822
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.343394
$ws.Range("N11").Value = 1.030182
$ws.Range("O11").Value = 0.5054542617721824
$ws.Range("P11").Value = 0.5054542617721824
$ws.Range("Q11").Value = 2.163265102646
$ws.Range("R11").Value = 19.469385923814
$ws.Range("S11").Value = 0.02907648638970841
$ws.Range("T11").Value = 0.02907648638970842
